# Update cryptocurrency price/volume data (refreshed snapshot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'75.218.41"
$ws.Cells.Item(2, 5).Value = "  +7.48%  "
$ws.Cells.Item(3, 4).Value = "'2.664.45"
$ws.Cells.Item(3, 5).Value = "  +8.94%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'186.83"
$ws.Cells.Item(5, 5).Value = "  +11.89%  "
$ws.Cells.Item(6, 4).Value = "'586.32"
$ws.Cells.Item(6, 5).Value = "  +3.31%  "
$ws.Cells.Item(7, 5).Value = "  -0.14%  "
$ws.Cells.Item(8, 4).Value = "'0.533"
$ws.Cells.Item(8, 5).Value = "  +3.73%  "
$ws.Cells.Item(9, 5).Value = "  +11.21%  "
$ws.Cells.Item(10, 4).Value = "'2.663.82"
$ws.Cells.Item(10, 5).Value = "  +8.97%  "
$ws.Cells.Item(11, 5).Value = "  +1.35%  "
$ws.Cells.Item(12, 5).Value = "  +6.18%  "
$ws.Cells.Item(13, 4).Value = "'4.71"
$ws.Cells.Item(13, 5).Value = "  +0.05%  "
$ws.Cells.Item(14, 4).Value = "'75.007.14"
$ws.Cells.Item(14, 5).Value = "  +7.13%  "
$ws.Cells.Item(15, 4).Value = "'3.156.50"
$ws.Cells.Item(15, 5).Value = "  +9.02%  "
$ws.Cells.Item(16, 5).Value = "  +3.80%  "
$ws.Cells.Item(17, 4).Value = "'26.48"
$ws.Cells.Item(17, 5).Value = "  +9.90%  "
$ws.Cells.Item(18, 4).Value = "'2.698.68"
$ws.Cells.Item(18, 5).Value = "  +10.40%  "
$ws.Cells.Item(19, 4).Value = "'9.18"
$ws.Cells.Item(19, 5).Value = "  +28.69%  "
$ws.Cells.Item(20, 5).Value = "  +10.09%  "
$ws.Cells.Item(21, 4).Value = "'370.58"
$ws.Cells.Item(21, 5).Value = "  +8.61%  "
$ws.Cells.Item(22, 4).Value = "'2.27"
$ws.Cells.Item(22, 5).Value = "  +13.65%  "
$ws.Cells.Item(23, 4).Value = "'4.06"
$ws.Cells.Item(23, 5).Value = "  +4.39%  "
$ws.Cells.Item(24, 5).Value = "  +3.64%  "
$ws.Cells.Item(25, 4).Value = "'0.999"
$ws.Cells.Item(25, 5).Value = "  -0.03%  "
$ws.Cells.Item(26, 4).Value = "'69.62"
$ws.Cells.Item(26, 5).Value = "  +5.05%  "
$ws.Cells.Item(27, 4).Value = "'4.14"
$ws.Cells.Item(27, 5).Value = "  +8.47%  "
$ws.Cells.Item(28, 4).Value = "'9.29"
$ws.Cells.Item(28, 5).Value = "  +9.54%  "
$ws.Cells.Item(29, 4).Value = "'2.795.12"
$ws.Cells.Item(29, 5).Value = "  +8.65%  "
$ws.Cells.Item(30, 5).Value = "  +0.59%  "
$ws.Cells.Item(31, 4).Value = "'0.0₃0942"
$ws.Cells.Item(31, 5).Value = "  +10.06%  "
$ws.Cells.Item(32, 5).Value = "  +14.39%  "
$ws.Cells.Item(33, 4).Value = "'518.64"
$ws.Cells.Item(33, 5).Value = "  +13.48%  "
$ws.Cells.Item(34, 5).Value = "  +3.84%  "
$ws.Cells.Item(35, 5).Value = "  +7.72%  "
$ws.Cells.Item(36, 4).Value = "'0.999"
$ws.Cells.Item(36, 5).Value = "  +0.00%  "
$ws.Cells.Item(37, 4).Value = "'163.44"
$ws.Cells.Item(37, 5).Value = "  +2.63%  "
$ws.Cells.Item(38, 5).Value = "  +6.06%  "
$ws.Cells.Item(39, 4).Value = "'19.13"
$ws.Cells.Item(39, 5).Value = "  +5.02%  "
$ws.Cells.Item(40, 4).Value = "'19.36"
$ws.Cells.Item(40, 5).Value = "  +1.46%  "
$ws.Cells.Item(41, 5).Value = "  +0.04%  "
$ws.Cells.Item(42, 4).Value = "'4.97"
$ws.Cells.Item(42, 5).Value = "  +12.75%  "
$ws.Cells.Item(43, 4).Value = "'168.93"
$ws.Cells.Item(43, 5).Value = "  +25.30%  "
$ws.Cells.Item(44, 4).Value = "'1.70"
$ws.Cells.Item(44, 5).Value = "  +11.29%  "
$ws.Cells.Item(45, 4).Value = "'0.328"
$ws.Cells.Item(45, 5).Value = "  +8.24%  "
$ws.Cells.Item(46, 5).Value = "  +9.78%  "
$ws.Cells.Item(47, 4).Value = "'2.36"
$ws.Cells.Item(47, 5).Value = "  +11.27%  "
$ws.Cells.Item(48, 4).Value = "'39.10"
$ws.Cells.Item(48, 5).Value = "  +2.83%  "
$ws.Cells.Item(49, 5).Value = "  +16.26%  "
$ws.Cells.Item(50, 5).Value = "  +7.22%  "
$ws.Cells.Item(51, 4).Value = "'0.533"
$ws.Cells.Item(51, 5).Value = "  +8.67%  "
